$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new time registration row for Thomas Borg
$ws.Range("A10").Value = (Get-Date -Year 2015 -Month 3 -Day 6).Date
$ws.Range("D10").Value = "Udarbejdelse af AD + Dataordbog for UC 6"
$ws.Range("B10").Value = "5 hrs?"
$ws.Range("C10").Value = "Requirement Specifier"

$ws.Range("B11").Select()
